# Aufgabe 7 - refactored evaluation (grading of Aufgabe 6 worksheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Evaluation")

# B18: change dropdown answer from "yes" to "no"
$ws.Range("B18").Value = "no"

# B19: change dropdown answer from "yes" to "no" and highlight it
# in red-on-green to flag it for special feedback attention.
# Copy B18's existing (light-green) cell format first so we reuse the
# existing fill, then switch the font color to red.
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B19").Value = "no"
$ws.Range("B19").Font.Color = 255

# Add feedback references in column H next to the affected criteria
$ws.Range("H18").Value = "siehe Feedback-Datei "
$ws.Range("H19").Value = "siehe Feedback-Datei "
$ws.Range("H45").Value = "siehe Feedback, könnte man geschickter machen, zu viele Ausgaben"

# Update the last active selection on the sheet
$ws.Range("E49").Select()

$wb.Save()
